# Auto-generated Excel COM-interop script to apply scheduled market-data update
# to the Ultima_Profits workbook (columns H-N: price/profit figures per leve row).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 147.2
$ws.Range("I5").Value = 158.5
$ws.Range("J5").Value = 102
$ws.Range("K5").Value = 158.5
$ws.Range("L5").Value = 102
$ws.Range("M5").Value = -43.5
$ws.Range("N5").Value = -332

$ws.Range("H43").Value = 3947.8823
$ws.Range("I43").Value = 1039.1875
$ws.Range("J43").Value = 6533.3887
$ws.Range("K43").Value = 1039.1875
$ws.Range("L43").Value = 6533.3887
$ws.Range("M43").Value = -970.1875
$ws.Range("N43").Value = -6671.3887

$ws.Range("H64").Value = 2265429.8
$ws.Range("I64").Value = 4051429.5
$ws.Range("J64").Value = 3163.2
$ws.Range("K64").Value = 4051429.5
$ws.Range("L64").Value = 3163.2
$ws.Range("M64").Value = -4051181.5
$ws.Range("N64").Value = -3659.2

$ws.Range("H67").Value = 2265429.8
$ws.Range("I67").Value = 4051429.5
$ws.Range("J67").Value = 3163.2
$ws.Range("K67").Value = 4051429.5
$ws.Range("L67").Value = 3163.2
$ws.Range("M67").Value = -4050571.5
$ws.Range("N67").Value = -4879.2

$ws.Range("H121").Value = 1257.5454
$ws.Range("J121").Value = 1414.7778
$ws.Range("L121").Value = 4244.3334
$ws.Range("N121").Value = -7738.3334

$ws.Range("H123").Value = 32926.668
$ws.Range("J123").Value = 32926.668
$ws.Range("L123").Value = 32926.668
$ws.Range("N123").Value = -42726.668

$ws.Range("H137").Value = 10527807
$ws.Range("I137").Value = 934.8
$ws.Range("J137").Value = 22224332
$ws.Range("K137").Value = 2804.4
$ws.Range("L137").Value = 66672996
$ws.Range("M137").Value = -254.3999999999996
$ws.Range("N137").Value = -66678096

$ws.Range("H138").Value = 2436.0732
$ws.Range("I138").Value = 1748.3889
$ws.Range("J138").Value = 2974.261
$ws.Range("K138").Value = 5245.1667
$ws.Range("L138").Value = 8922.782999999999
$ws.Range("M138").Value = -105.1666999999998
$ws.Range("N138").Value = -19202.783


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 687.73334
$ws.Range("I2").Value = 445.72726
$ws.Range("J2").Value = 1353.25
$ws.Range("K2").Value = 445.72726
$ws.Range("L2").Value = 1353.25
$ws.Range("M2").Value = -332.72726
$ws.Range("N2").Value = -1579.25

$ws.Range("H32").Value = 12983.821
$ws.Range("I32").Value = 11814.458
$ws.Range("J32").Value = 20000
$ws.Range("K32").Value = 11814.458
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = -11527.458
$ws.Range("N32").Value = -20574

$ws.Range("H45").Value = 1484.2759
$ws.Range("I45").Value = 1260
$ws.Range("J45").Value = 1910.4
$ws.Range("K45").Value = 1260
$ws.Range("L45").Value = 1910.4
$ws.Range("M45").Value = -883
$ws.Range("N45").Value = -2664.4

$ws.Range("H63").Value = 62500988
$ws.Range("I63").Value = 83334370
$ws.Range("J63").Value = 853
$ws.Range("K63").Value = 83334370
$ws.Range("L63").Value = 853
$ws.Range("M63").Value = -83333684
$ws.Range("N63").Value = -2225

$ws.Range("H66").Value = 62500988
$ws.Range("I66").Value = 83334370
$ws.Range("J66").Value = 853
$ws.Range("K66").Value = 416671850
$ws.Range("L66").Value = 4265
$ws.Range("M66").Value = -416668418
$ws.Range("N66").Value = -11129

$ws.Range("H80").Value = 24483.777
$ws.Range("J80").Value = 24483.777
$ws.Range("L80").Value = 24483.777
$ws.Range("N80").Value = -26479.777

$ws.Range("H83").Value = 24483.777
$ws.Range("J83").Value = 24483.777
$ws.Range("L83").Value = 73451.33099999999
$ws.Range("N83").Value = -83435.33099999999

$ws.Range("H97").Value = 4442.269
$ws.Range("I97").Value = 5041.773
$ws.Range("J97").Value = 1145
$ws.Range("K97").Value = 5041.773
$ws.Range("L97").Value = 1145
$ws.Range("M97").Value = -4545.773
$ws.Range("N97").Value = -2137

$ws.Range("H116").Value = 687.73334
$ws.Range("I116").Value = 445.72726
$ws.Range("J116").Value = 1353.25
$ws.Range("K116").Value = 445.72726
$ws.Range("L116").Value = 1353.25
$ws.Range("M116").Value = 1848.27274
$ws.Range("N116").Value = -5941.25


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 687.73334
$ws.Range("I3").Value = 445.72726
$ws.Range("J3").Value = 1353.25
$ws.Range("K3").Value = 445.72726
$ws.Range("L3").Value = 1353.25
$ws.Range("M3").Value = -331.72726
$ws.Range("N3").Value = -1581.25


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5379482
$ws.Range("I31").Value = 3429.0435
$ws.Range("J31").Value = 20835634
$ws.Range("K31").Value = 3429.0435
$ws.Range("L31").Value = 20835634
$ws.Range("M31").Value = -3134.0435
$ws.Range("N31").Value = -20836224

$ws.Range("H34").Value = 5379482
$ws.Range("I34").Value = 3429.0435
$ws.Range("J34").Value = 20835634
$ws.Range("K34").Value = 3429.0435
$ws.Range("L34").Value = 20835634
$ws.Range("M34").Value = -3227.0435
$ws.Range("N34").Value = -20836038

$ws.Range("H58").Value = 1811
$ws.Range("I58").Value = 699.1905
$ws.Range("J58").Value = 5702.3335
$ws.Range("K58").Value = 699.1905
$ws.Range("L58").Value = 5702.3335
$ws.Range("M58").Value = -496.1905
$ws.Range("N58").Value = -6108.3335

$ws.Range("H62").Value = 2686.1904
$ws.Range("I62").Value = 2207.6924
$ws.Range("J62").Value = 3463.75
$ws.Range("K62").Value = 2207.6924
$ws.Range("L62").Value = 3463.75
$ws.Range("M62").Value = -1583.6924
$ws.Range("N62").Value = -4711.75

$ws.Range("H65").Value = 2686.1904
$ws.Range("I65").Value = 2207.6924
$ws.Range("J65").Value = 3463.75
$ws.Range("K65").Value = 11038.462
$ws.Range("L65").Value = 17318.75
$ws.Range("M65").Value = -7918.462
$ws.Range("N65").Value = -23558.75

$ws.Range("H122").Value = 1898.1765
$ws.Range("I122").Value = 1898.1765
$ws.Range("K122").Value = 5694.529500000001
$ws.Range("M122").Value = -3244.529500000001

$ws.Range("H129").Value = 43333
$ws.Range("J129").Value = 43333
$ws.Range("L129").Value = 43333
$ws.Range("N129").Value = -53333

$ws.Range("H132").Value = 14707814
$ws.Range("I132").Value = 21740810
$ws.Range("J132").Value = 2457.9092
$ws.Range("K132").Value = 65222430
$ws.Range("L132").Value = 7373.7276
$ws.Range("M132").Value = -65219900
$ws.Range("N132").Value = -12433.7276

$ws.Range("H136").Value = 1811
$ws.Range("I136").Value = 699.1905
$ws.Range("J136").Value = 5702.3335
$ws.Range("K136").Value = 2097.5715
$ws.Range("L136").Value = 17107.0005
$ws.Range("M136").Value = 452.4285
$ws.Range("N136").Value = -22207.0005


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 885.125
$ws.Range("I68").Value = 1014.8571
$ws.Range("J68").Value = 784.2222
$ws.Range("K68").Value = 3044.5713
$ws.Range("L68").Value = 2352.6666
$ws.Range("M68").Value = -2233.5713
$ws.Range("N68").Value = -3974.6666

$ws.Range("H71").Value = 885.125
$ws.Range("I71").Value = 1014.8571
$ws.Range("J71").Value = 784.2222
$ws.Range("K71").Value = 9133.713899999999
$ws.Range("L71").Value = 7057.999800000001
$ws.Range("M71").Value = -5077.713899999999
$ws.Range("N71").Value = -15169.9998

$ws.Range("H98").Value = 837.8333
$ws.Range("I98").Value = 710.8333
$ws.Range("J98").Value = 901.3333
$ws.Range("K98").Value = 2132.4999
$ws.Range("L98").Value = 2703.9999
$ws.Range("M98").Value = -634.4998999999998
$ws.Range("N98").Value = -5699.9999

$ws.Range("H107").Value = 942.1053000000001
$ws.Range("I107").Value = 159.71428
$ws.Range("J107").Value = 1398.5
$ws.Range("K107").Value = 479.14284
$ws.Range("L107").Value = 4195.5
$ws.Range("M107").Value = 1440.85716
$ws.Range("N107").Value = -8035.5

$ws.Range("H131").Value = 1233
$ws.Range("I131").Value = 520
$ws.Range("J131").Value = 1411.25
$ws.Range("K131").Value = 1560
$ws.Range("L131").Value = 4233.75
$ws.Range("M131").Value = 3480
$ws.Range("N131").Value = -14313.75


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 8458
$ws.Range("I57").Value = 3200
$ws.Range("J57").Value = 15030.5
$ws.Range("K57").Value = 3200
$ws.Range("L57").Value = 15030.5
$ws.Range("M57").Value = -2380
$ws.Range("N57").Value = -16670.5

$ws.Range("H70").Value = 10836.434
$ws.Range("I70").Value = 25247
$ws.Range("J70").Value = 4660.476
$ws.Range("K70").Value = 25247
$ws.Range("L70").Value = 4660.476
$ws.Range("M70").Value = -24977
$ws.Range("N70").Value = -5200.476

$ws.Range("H73").Value = 10836.434
$ws.Range("I73").Value = 25247
$ws.Range("J73").Value = 4660.476
$ws.Range("K73").Value = 25247
$ws.Range("L73").Value = 4660.476
$ws.Range("M73").Value = -24311
$ws.Range("N73").Value = -6532.476

$ws.Range("H97").Value = 1714.55
$ws.Range("I97").Value = 1604
$ws.Range("J97").Value = 2046.2
$ws.Range("K97").Value = 1604
$ws.Range("L97").Value = 2046.2
$ws.Range("M97").Value = -1108
$ws.Range("N97").Value = -3038.2

$ws.Range("H123").Value = 44375
$ws.Range("J123").Value = 44375
$ws.Range("L123").Value = 44375
$ws.Range("N123").Value = -49275


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 9284.615
$ws.Range("I122").Value = 19333.334
$ws.Range("J122").Value = 6270
$ws.Range("K122").Value = 58000.00199999999
$ws.Range("L122").Value = 18810
$ws.Range("M122").Value = -55550.00199999999
$ws.Range("N122").Value = -23710


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4451.8887
$ws.Range("I62").Value = 3983.3333
$ws.Range("J62").Value = 4826.7334
$ws.Range("K62").Value = 3983.3333
$ws.Range("L62").Value = 4826.7334
$ws.Range("M62").Value = -3359.3333
$ws.Range("N62").Value = -6074.7334

$ws.Range("H65").Value = 4451.8887
$ws.Range("I65").Value = 3983.3333
$ws.Range("J65").Value = 4826.7334
$ws.Range("K65").Value = 19916.6665
$ws.Range("L65").Value = 24133.667
$ws.Range("M65").Value = -16796.6665
$ws.Range("N65").Value = -30373.667

